$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = @(
    @('Website Name','Profile URL'),
    @('35photo','https://35photo.pro/@gorbash/'),
    @('Bandlab','https://www.bandlab.com/api/v1.3/users/gorbash'),
    @('Audiojungle','https://audiojungle.net/user/gorbash'),
    @('BodyBuilding.com','http://api.bodybuilding.com/api-proxy/bbc/get?slug=gorbash'),
    @('Chess.com','https://api.chess.com/pub/player/gorbash'),
    @('chaturbate','https://chaturbate.com/gorbash/'),
    @('cfx.re','https://forum.cfx.re/u/gorbash.json'),
    @('chatango.com','https://gorbash.chatango.com'),
    @('cda.pl','https://www.cda.pl/gorbash'),
    @('Chyoa','https://chyoa.com/user/gorbash'),
    @('Chomikuj.pl','https://chomikuj.pl/gorbash/'),
    @('Code Project','https://www.codeproject.com/Members/gorbash'),
    @('Codewars','https://www.codewars.com/users/gorbash'),
    @('cHEEZburger','https://profile.cheezburger.com/gorbash'),
    @('crowdin','https://crowdin.com/profile/gorbash'),
    @('Mastodon-Defcon','https://defcon.social/@gorbash'),
    @('diigo','https://www.diigo.com/interact_api/load_profile_info?name=gorbash'),
    @('Disqus','https://disqus.com/by/gorbash/'),
    @('DeviantArt','https://www.deviantart.com/gorbash'),
    @('Discogs','https://api.discogs.com/users/gorbash'),
    @('DockerHub','https://hub.docker.com/v2/users/gorbash/'),
    @('Dating.ru','https://dating.ru/gorbash/'),
    @('datezone','https://www.datezone.com/users/gorbash/'),
    @('Duolingo','https://www.duolingo.com/2017-06-30/users?username=gorbash&_=1628308619574'),
    @('Cults3D','https://cults3d.com/en/users/gorbash/creations'),
    @('Etsy','https://www.etsy.com/people/gorbash'),
    @('Demotywatory','https://demotywatory.pl/user/gorbash'),
    @('eBay','https://www.ebay.com/usr/gorbash'),
    @('Fark','https://www.fark.com/users/gorbash'),
    @('FatSecret','https://www.fatsecret.com/member/gorbash'),
    @('Fabswingers','https://www.fabswingers.com/profile/gorbash'),
    @('Filmweb','https://www.filmweb.pl/user/gorbash'),
    @('fansly','https://apiv2.fansly.com/api/v1/account?usernames=gorbash'),
    @('Foursquare','https://foursquare.com/gorbash'),
    @('FurAffinity','https://www.furaffinity.net/user/gorbash'),
    @('Geocaching','https://www.geocaching.com/p/?u=gorbash'),
    @('Gravatar','https://en.gravatar.com/gorbash.json'),
    @('GitHub','https://github.com/gorbash'),
    @('giters','https://giters.com/gorbash'),
    @('GitLab','https://gitlab.com/gorbash'),
    @('HackerOne','https://hackerone.com/gorbash'),
    @('Hacker News','https://news.ycombinator.com/user?id=gorbash'),
    @('HudsonRock','https://cavalier.hudsonrock.com/api/json/v2/osint-tools/search-by-username?username=gorbash'),
    @('IFTTT','https://ifttt.com/p/gorbash'),
    @('FriendFinder-X','https://www.friendfinder-x.com/profile/gorbash'),
    @('hugging_face','https://huggingface.co/gorbash'),
    @('imgur','https://api.imgur.com/account/v1/accounts/gorbash?client_id=546c25a59c58ad7&include=trophies%2Cmedallions'),
    @('inaturalist','https://inaturalist.nz/people/gorbash'),
    @('InkBunny','https://inkbunny.net/gorbash'),
    @('issuu','https://issuu.com/gorbash'),
    @('Instagram2','https://dumpoir.com/v/gorbash'),
    @('kaggle','https://www.kaggle.com/gorbash'),
    @('Jeuxvideo','https://www.jeuxvideo.com/profil/gorbash?mode=infos'),
    @('Internet Archive User Search','https://archive.org/advancedsearch.php?q=gorbash&output=json'),
    @('Keybase','https://keybase.io/gorbash'),
    @('Kongregate','https://www.kongregate.com/accounts/gorbash'),
    @('Linktree','https://linktr.ee/gorbash'),
    @('ask.fm','https://ask.fm/gorbash'),
    @('waytohey','https://waytohey.com/gorbash'),
    @('LibraryThing','https://www.librarything.com/profile/gorbash'),
    @('Livejournal','https://gorbash.livejournal.com'),
    @('Mastodon API','https://mastodon.social/api/v2/search?q=gorbash&limit=1&type=accounts'),
    @('MCUUID (Minecraft)','https://playerdb.co/api/player/minecraft/gorbash'),
    @('Mod DB','https://www.moddb.com/members/gorbash'),
    @('MySpace','https://myspace.com/gorbash'),
    @('MyAnimeList','https://myanimelist.net/profile/gorbash'),
    @('OpenStreetMap','https://www.openstreetmap.org/user/gorbash'),
    @('MCName (Minecraft)','https://mcname.info/en/search?q=gorbash'),
    @('palnet','https://www.palnet.io/@gorbash/'),
    @('Pokerstrategy','http://www.pokerstrategy.net/user/gorbash/profile/'),
    @('Playstation Network','https://psnprofiles.com/xhr/search/users?q=gorbash'),
    @('Periscope','https://www.periscope.tv/gorbash'),
    @('pikabu','https://pikabu.ru/@gorbash'),
    @('Pinterest','https://www.pinterest.com/gorbash/'),
    @('Flipboard','https://flipboard.com/@gorbash'),
    @('prv.pl','https://www.prv.pl/osoba/gorbash'),
    @('Pornhub Users','https://www.pornhub.com/users/gorbash'),
    @('public','https://public.com/@gorbash'),
    @('Roblox','https://auth.roblox.com/v1/usernames/validate?username=gorbash&birthday=2019-12-31T23:00:00.000Z'),
    @('rsi','https://robertsspaceindustries.com/citizens/gorbash'),
    @('sentimente','https://www.sentimente.com/amp/gorbash.html'),
    @('RumbleUser','https://rumble.com/user/gorbash'),
    @('SFD','https://www.sfd.pl/profile/gorbash'),
    @('slideshare','https://www.slideshare.net/gorbash'),
    @('SoundCloud','https://soundcloud.com/gorbash'),
    @('smule','https://www.smule.com/api/profile/?handle=gorbash'),
    @('Sourceforge','https://sourceforge.net/u/gorbash/profile'),
    @('sofurry','https://gorbash.sofurry.com'),
    @('Snapchat Stories','https://story.snapchat.com/s/gorbash'),
    @('Spotify','https://open.spotify.com/user/gorbash'),
    @('Steam','https://steamcommunity.com/id/gorbash'),
    @('Telegram','https://t.me/gorbash'),
    @('tabletoptournament','https://www.tabletoptournaments.net/eu/player/gorbash'),
    @('Trello','https://trello.com/1/Members/gorbash?fields=activityBlocked%2CavatarUrl%2Cbio%2CbioData%2Cconfirmed%2CfullName%2CidEnterprise%2CidMemberReferrer%2Cinitials%2CmemberType%2CnonPublic%2Cproducts%2Curl%2Cusername'),
    @('themeforest','https://themeforest.net/user/gorbash'),
    @('TikTok','https://www.tiktok.com/oembed?url=https://www.tiktok.com/@gorbash'),
    @('trakt','https://trakt.tv/users/gorbash'),
    @('tumblr','https://gorbash.tumblr.com'),
    @('Twitter','https://nitter.privacydev.net/gorbash'),
    @('Twitch','https://twitchtracker.com/gorbash'),
    @('Twitter archived tweets','http://archive.org/wayback/available?url=https://twitter.com/gorbash/status/*'),
    @('Vimeo','https://vimeo.com/gorbash'),
    @('Venmo','https://account.venmo.com/u/gorbash'),
    @('untappd','https://untappd.com/user/gorbash/'),
    @('wattpad','https://www.wattpad.com/api/v3/users/gorbash?fields=username%2Cname%2Cdescription%2Cavatar%2CbackgroundUrl%2CcreateDate%2Clocation%2Cfollowing%2CfollowingRequest%2CnumFollowing%2Cfollower%2CfollowerRequest%2CnumFollowers%2CnumLists%2CnumStoriesPublished%2CvotesReceived%2Cfacebook%2Ctwitter%2Cwebsite%2Csmashwords%2Chighlight_colour%2Chtml_enabled%2Cverified%2Cambassador%2Cwattpad_squad%2Cis_staff%2Cprograms(wattpad_stars)%2CisPrivate%2CisMuted%2CexternalId%2Cnotes'),
    @('Wikidot','http://www.wikidot.com/user:info/gorbash'),
    @('VK','https://vk.com/gorbash'),
    @('Xbox Gamertag','https://www.xboxgamertag.com/search/gorbash'),
    @('WordPress Support','https://wordpress.org/support/users/gorbash/'),
    @('xHamster','https://xhamster.com/users/gorbash'),
    @('Zbiornik','https://mini.zbiornik.com/gorbash'),
    @('XVIDEOS-profiles','https://www.xvideos.com/profiles/gorbash'),
    @('XVIDEOS-models','https://www.xvideos.com/models/gorbash'),
    @('YouTube User','https://www.youtube.com/user/gorbash/about'),
    @('YouTube User2','https://www.youtube.com/@gorbash'),
    @('XNXX','https://www.xnxx.com/mobile/profile/gorbash')
)
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($i+1, 1).Value = $row[0]
    $ws.Cells.Item($i+1, 2).Value = $row[1]
}
